# Applies the "Updated cryptos list" data refresh: column D (Price) and
# column E (Volume(1h)) values are updated for rows 2-51 of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values look like plain numbers (e.g. "1.00",
# "317.86"). The source cells are text (inline strings), so before writing
# those values we temporarily force the cells number format to Text ("@")
# to stop Excel from auto-converting the assigned value into a numeric cell.
$textForcedCells = @(
    "D4",
    "D5",
    "D6",
    "D10",
    "D14",
    "D15",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D27",
    "D29",
    "D30",
    "D31",
    "D34",
    "D35",
    "D36",
    "D39",
    "D41",
    "D44",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Column D - Price
$ws.Range("D2").Value = "41.633.00"
$ws.Range("D3").Value = "2.472.71"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "317.86"
$ws.Range("D6").Value = "92.55"
$ws.Range("D10").Value = "0.0870"
$ws.Range("D13").Value = "2.853.40"
$ws.Range("D14").Value = "6.90"
$ws.Range("D15").Value = "15.57"
$ws.Range("D16").Value = "2.473.10"
$ws.Range("D17").Value = "0.788"
$ws.Range("D18").Value = "41.590.33"
$ws.Range("D19").Value = "0.0₃0955"
$ws.Range("D20").Value = "6.49"
$ws.Range("D21").Value = "71.28"
$ws.Range("D22").Value = "11.35"
$ws.Range("D23").Value = "240.52"
$ws.Range("D27").Value = "24.68"
$ws.Range("D29").Value = "9.92"
$ws.Range("D30").Value = "36.39"
$ws.Range("D31").Value = "158.63"
$ws.Range("D34").Value = "0.0771"
$ws.Range("D35").Value = "2.58"
$ws.Range("D36").Value = "17.38"
$ws.Range("D39").Value = "0.116"
$ws.Range("D41").Value = "3.99"
$ws.Range("D43").Value = "1.984.90"
$ws.Range("D44").Value = "19.18"
$ws.Range("D46").Value = "3.01"
$ws.Range("D47").Value = "9.24"
$ws.Range("D48").Value = "2.710.58"
$ws.Range("D49").Value = "97.50"
$ws.Range("D50").Value = "67.32"
$ws.Range("D51").Value = "73.47"

# Put the number format back to the default (General/Normal style) now
# that the text values have been written, so no stray style stays applied.
foreach ($cellRef in $textForcedCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Column E - Volume(1h)
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  +1.47%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("E10").Value = "  +10.69%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  +1.79%  "
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  +3.75%  "
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("E32").Value = "  +1.30%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  -2.51%  "
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("E45").Value = "  +0.86%  "
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("E51").Value = "  +1.26%  "

